$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the confidential "as of" date string (A10)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-12 for illustrative purposes only and are subject to change."

# Update numeric holdings values (columns D and E, rows 2-7)
$ws.Range("D2").Value = 0.4775578854398405
$ws.Range("E2").Value = -0.01138594424813488

$ws.Range("D3").Value = 0.3409259820063416
$ws.Range("E3").Value = -0.01533396048918145

$ws.Range("D4").Value = 0.09605334816871629
$ws.Range("E4").Value = -0.02463142754404879

$ws.Range("D5").Value = 0.05368369178998476
$ws.Range("E5").Value = -0.004825367647058876

$ws.Range("D6").Value = 0.03177909259511678
$ws.Range("E6").Value = -0.006113537117903856

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = -0.01348445029416845
